$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.469240333333333
$ws.Range("H2").Value = 13.407721
$ws.Range("I2").Value = 0.3358648218165975
$ws.Range("J2").Value = 0.3358648218165975
$ws.Range("M2").Value = 1.819857
$ws.Range("N2").Value = 5.459571
$ws.Range("O2").Value = 0.01485317462584607
$ws.Range("P2").Value = 0.01485317462584607
$ws.Range("Q2").Value = 8.133378305298999
$ws.Range("R2").Value = 73.200404747691
$ws.Range("S2").Value = 0.004988658849120597
$ws.Range("T2").Value = 0.004988658849120598
$ws.Range("G3").Value = 4.469240333333333
$ws.Range("H3").Value = 13.407721
$ws.Range("I3").Value = 0.3358648218165975
$ws.Range("J3").Value = 0.3358648218165975
$ws.Range("O3").Value = 0.726618572334523
$ws.Range("P3").Value = 0.7266185723345231
$ws.Range("Q3").Value = 397.8855619302532
$ws.Range("R3").Value = 3580.970057372279
$ws.Range("S3").Value = 0.244045617325765
$ws.Range("T3").Value = 0.2440456173257651
$ws.Range("G4").Value = 4.469240333333333
$ws.Range("H4").Value = 13.407721
$ws.Range("I4").Value = 0.3358648218165975
$ws.Range("J4").Value = 0.3358648218165975
$ws.Range("M4").Value = 31.52924033333333
$ws.Range("N4").Value = 94.58772099999999
$ws.Range("O4").Value = 0.257333028084772
$ws.Range("P4").Value = 0.257333028084772
$ws.Range("Q4").Value = 140.9117525770934
$ws.Range("R4").Value = 1268.205773193841
$ws.Range("S4").Value = 0.08642911162521745
$ws.Range("T4").Value = 0.08642911162521745
$ws.Range("G5").Value = 4.469240333333333
$ws.Range("H5").Value = 13.407721
$ws.Range("I5").Value = 0.3358648218165975
$ws.Range("J5").Value = 0.3358648218165975
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1464426666666667
$ws.Range("N5").Value = 0.439328
$ws.Range("O5").Value = 0.001195224954858853
$ws.Range("P5").Value = 0.001195224954858853
$ws.Range("Q5").Value = 0.6544874723875554
$ws.Range("R5").Value = 5.890387251488
$ws.Range("S5").Value = 0.0004014340164944193
$ws.Range("T5").Value = 0.0004014340164944194
$ws.Range("I6").Value = 0.2721973992379558
$ws.Range("J6").Value = 0.2721973992379558
$ws.Range("M6").Value = 1.819857
$ws.Range("N6").Value = 5.459571
$ws.Range("O6").Value = 0.01485317462584607
$ws.Range("P6").Value = 0.01485317462584607
$ws.Range("Q6").Value = 6.591593635042
$ws.Range("R6").Value = 59.324342715378
$ws.Range("S6").Value = 0.004042995503582497
$ws.Range("T6").Value = 0.004042995503582497
$ws.Range("I7").Value = 0.2721973992379558
$ws.Range("J7").Value = 0.2721973992379558
$ws.Range("O7").Value = 0.726618572334523
$ws.Range("P7").Value = 0.7266185723345231
$ws.Range("S7").Value = 0.1977836856274536
$ws.Range("T7").Value = 0.1977836856274536
$ws.Range("I8").Value = 0.2721973992379558
$ws.Range("J8").Value = 0.2721973992379558
$ws.Range("M8").Value = 31.52924033333333
$ws.Range("N8").Value = 94.58772099999999
$ws.Range("O8").Value = 0.257333028084772
$ws.Range("P8").Value = 0.257333028084772
$ws.Range("Q8").Value = 114.2001486374531
$ws.Range("R8").Value = 1027.801337737078
$ws.Range("S8").Value = 0.0700453809827028
$ws.Range("T8").Value = 0.07004538098270278
$ws.Range("I9").Value = 0.2721973992379558
$ws.Range("J9").Value = 0.2721973992379558
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1464426666666667
$ws.Range("N9").Value = 0.439328
$ws.Range("O9").Value = 0.001195224954858853
$ws.Range("P9").Value = 0.001195224954858853
$ws.Range("Q9").Value = 0.5304210987448889
$ws.Range("R9").Value = 4.773789888704
$ws.Range("S9").Value = 0.0003253371242168828
$ws.Range("T9").Value = 0.0003253371242168828
$ws.Range("G10").Value = 3.887787333333333
$ws.Range("H10").Value = 11.663362
$ws.Range("I10").Value = 0.2921684453243378
$ws.Range("J10").Value = 0.2921684453243377
$ws.Range("M10").Value = 1.819857
$ws.Range("N10").Value = 5.459571
$ws.Range("O10").Value = 0.01485317462584607
$ws.Range("P10").Value = 0.01485317462584607
$ws.Range("Q10").Value = 7.075216993078
$ws.Range("R10").Value = 63.676952937702
$ws.Range("S10").Value = 0.004339628938564348
$ws.Range("T10").Value = 0.004339628938564348
$ws.Range("G11").Value = 3.887787333333333
$ws.Range("H11").Value = 11.663362
$ws.Range("I11").Value = 0.2921684453243378
$ws.Range("J11").Value = 0.2921684453243377
$ws.Range("O11").Value = 0.726618572334523
$ws.Range("P11").Value = 0.7266185723345231
$ws.Range("Q11").Value = 346.1202200855733
$ws.Range("R11").Value = 3115.08198077016
$ws.Range("S11").Value = 0.2122950186227675
$ws.Range("T11").Value = 0.2122950186227675
$ws.Range("G12").Value = 3.887787333333333
$ws.Range("H12").Value = 11.663362
$ws.Range("I12").Value = 0.2921684453243378
$ws.Range("J12").Value = 0.2921684453243377
$ws.Range("M12").Value = 31.52924033333333
$ws.Range("N12").Value = 94.58772099999999
$ws.Range("O12").Value = 0.257333028084772
$ws.Range("P12").Value = 0.257333028084772
$ws.Range("Q12").Value = 122.5789811975558
$ws.Range("R12").Value = 1103.210830778002
$ws.Range("S12").Value = 0.075184590746132
$ws.Range("T12").Value = 0.07518459074613199
$ws.Range("G13").Value = 3.887787333333333
$ws.Range("H13").Value = 11.663362
$ws.Range("I13").Value = 0.2921684453243378
$ws.Range("J13").Value = 0.2921684453243377
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.1464426666666667
$ws.Range("N13").Value = 0.439328
$ws.Range("O13").Value = 0.001195224954858853
$ws.Range("P13").Value = 0.001195224954858853
$ws.Range("Q13").Value = 0.5693379445262222
$ws.Range("R13").Value = 5.124041500735999
$ws.Range("S13").Value = 0.0003492070168739628
$ws.Range("T13").Value = 0.0003492070168739628
$ws.Range("G14").Value = 1.327597
$ws.Range("H14").Value = 3.982791
$ws.Range("I14").Value = 0.09976933362110896
$ws.Range("J14").Value = 0.09976933362110893
$ws.Range("M14").Value = 1.819857
$ws.Range("N14").Value = 5.459571
$ws.Range("O14").Value = 0.01485317462584607
$ws.Range("P14").Value = 0.01485317462584607
$ws.Range("Q14").Value = 2.416036693629
$ws.Range("R14").Value = 21.744330242661
$ws.Range("S14").Value = 0.001481891334578626
$ws.Range("T14").Value = 0.001481891334578626
$ws.Range("G15").Value = 1.327597
$ws.Range("H15").Value = 3.982791
$ws.Range("I15").Value = 0.09976933362110896
$ws.Range("J15").Value = 0.09976933362110893
$ws.Range("O15").Value = 0.726618572334523
$ws.Range("P15").Value = 0.7266185723345231
$ws.Range("Q15").Value = 118.19272157332
$ws.Range("R15").Value = 1063.73449415988
$ws.Range("S15").Value = 0.07249425075853691
$ws.Range("T15").Value = 0.07249425075853691
$ws.Range("G16").Value = 1.327597
$ws.Range("H16").Value = 3.982791
$ws.Range("I16").Value = 0.09976933362110896
$ws.Range("J16").Value = 0.09976933362110893
$ws.Range("M16").Value = 31.52924033333333
$ws.Range("N16").Value = 94.58772099999999
$ws.Range("O16").Value = 0.257333028084772
$ws.Range("P16").Value = 0.257333028084772
$ws.Range("Q16").Value = 41.85812487881233
$ws.Range("R16").Value = 376.7231239093109
$ws.Range("S16").Value = 0.02567394473071982
$ws.Range("T16").Value = 0.02567394473071982
$ws.Range("G17").Value = 1.327597
$ws.Range("H17").Value = 3.982791
$ws.Range("I17").Value = 0.09976933362110896
$ws.Range("J17").Value = 0.09976933362110893
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.1464426666666667
$ws.Range("N17").Value = 0.439328
$ws.Range("O17").Value = 0.001195224954858853
$ws.Range("P17").Value = 0.001195224954858853
$ws.Range("Q17").Value = 0.1944168449386666
$ws.Range("R17").Value = 1.749751604448
$ws.Range("S17").Value = 0.0001192467972735878
$ws.Range("T17").Value = 0.0001192467972735877
